$d = $word.ActiveDocument

# Paragraph 1: title/author block - merge "John "+"Southworth" into one run,
# add a new "Keele University" line, and move the _GoBack bookmark here.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML(@'
<w:p w:rsidR="00EE4367" w:rsidRPr="00EE4367" w:rsidRDefault="00EE4367" w:rsidP="00EE4367"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00EE4367"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">High-precision photometry of Qatar-2 and WASP-55: the two currently known transiting planets in K2 field 6 </w:t></w:r><w:r w:rsidRPr="00EE4367"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:br/><w:t>John Southworth</w:t></w:r><w:r w:rsidRPr="00EE4367"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:br/><w:t>Keele University</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@)

# Paragraph 2: "The study of extrasolar planets..." - strip proofErr wrappers,
# merge all runs into a single run with identical formatting.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML(@'
<w:p w:rsidR="00EE4367" w:rsidRPr="00EE4367" w:rsidRDefault="00EE4367" w:rsidP="00EE4367"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00EE4367"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">The study of extrasolar planets is a young field but already represents one of the major scientific advances of the 21st century, due in particular to the remarkable results achieved during the main mission of the Kepler satellite. Of the known extrasolar planets, only those which transit their host stars can be characterised in detail. For these objects, analysis of photometry and spectroscopy can yield measurements of their masses and radii, and thus densities and surface gravities, enabling studies of the structure, formation and evolution of planets and planetary systems. </w:t></w:r></w:p>
'@)

# Paragraph 3: "The K2 mission field 6 ..." - strip proofErr wrappers,
# merge all runs into a single run with identical formatting.
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML(@'
<w:p w:rsidR="00EE4367" w:rsidRPr="00EE4367" w:rsidRDefault="00EE4367" w:rsidP="00EE4367"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00EE4367"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">The K2 mission field 6 contains two known transiting extrasolar planets: Qatar-2 and WASP-55. Both are gas-giant planets in short- period orbits (1.3 and 4.5 days respectively) whose masses and radii have been measured to precisions of between 2% and 10% from ground-based observations. Detailed error budgets from past studies of transiting planets have shown that the quality of the light curve is the main determinant of the precision of the measured properties of the system (e.g. Southworth, 2009, MNRAS, 394, 272), especially for density and surface gravity. </w:t></w:r></w:p>
'@)

# Paragraph 5: "We propose to observe Qatar-2 and WASP-55 ..." - strip proofErr
# wrappers, merge surrounding runs, and drop the _GoBack bookmark (moved to
# paragraph 1 above).
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertXML(@'
<w:p w:rsidR="00EE4367" w:rsidRPr="00EE4367" w:rsidRDefault="00EE4367" w:rsidP="00EE4367"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00EE4367"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>We propose to observe Qatar-2 and WASP-55 with Kepler K2 in short cadence. We will use the light curves to measure the physical properties of both systems to very high precision, using the methodology we have developed over the past 8 years (see Southworth, 2012, MNRAS, 426, 1291 and references therein) and a K2 data reduction pipeline currently being</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> developed by our group. Short-</w:t></w:r><w:r w:rsidRPr="00EE4367"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">cadence data is mandatory to resolve the spot crossing events in the transits of Qatar-2, which have a characteristic timescale of approximately 15 minutes. For WASP-55, short-cadence data is necessary to obtain a measurement of the planet's density to 2-3% precision (long cadence observations would only give a precision of 5-6%), which is needed to usefully constrain the internal structure of a planet of this mass (see Fortney et al., 2007, ApJ, 659, 1661). </w:t></w:r></w:p>
'@)

Write-Host "Done"
